$d = $word.ActiveDocument

# This index/TOC list currently starts with four bulleted paragraphs:
#   1. "Java pattern matching switch case statements. ..."
#   2. "Data Statements: Model"
#   3. "Schema Statements: Contexts"
#   4. "Interaction Statements: Views"
# They are being replaced by nine bulleted paragraphs with reworded /
# expanded text, dropping the incidental pageBreakBefore="0" and
# u val="none" paragraph formatting that the old paragraphs carried
# (the new paragraphs only need the list numbering + indent in pPr).

# Locate the block to replace by content rather than a hard-coded
# paragraph index, then expand each end to the full paragraph (so the
# trailing paragraph mark is included).
$rngStart = $d.Content.Duplicate
$rngStart.Find.Execute("Java pattern matching switch case statements", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rngStart.Expand(4) | Out-Null   # wdParagraph

$rngEnd = $d.Content.Duplicate
$rngEnd.Find.Execute("Interaction Statements: Views", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rngEnd.Expand(4) | Out-Null     # wdParagraph

$target = $d.Range($rngStart.Start, $rngEnd.End)

$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

function New-BulletParagraphXml([string]$text) {
    return "<w:p $wNs><w:pPr><w:numPr><w:ilvl w:val='0'/><w:numId w:val='2'/></w:numPr>" +
           "<w:ind w:left='600' w:hanging='360'/></w:pPr><w:r><w:rPr><w:rtl w:val='0'/></w:rPr>" +
           "<w:t xml:space='preserve'>$text</w:t></w:r></w:p>"
}

$newTexts = @(
    "Java pattern matching statements case classes. Resource Monad wrapping Case classes.",
    "C, S, P, O URN Case classes.",
    "SK, PK, OK Kinds Case classes.",
    "Data Statements: Model (Kinds Aggregation) Case classes.",
    "Schema Statements: Contexts (Kinds Context Alignment) Case classes.",
    "Interaction Statements: Views (Context instances Activation) Case clases.",
    "Composition of Case classes instances / Augmentations via pattern matching (Aggregation, Alignment, Activation).",
    "Functors / Transforms: over Resource wrapping Case classes.",
    "Functors / Transforms: over (reified) Case classes instances higher kinds."
)

$bulletsXml = ($newTexts | ForEach-Object { New-BulletParagraphXml $_ }) -join ""

# InsertXML on a range that spans whole paragraphs (including their end
# marks) replaces those paragraphs outright, letting us control pPr/rPr
# precisely instead of inheriting the old pageBreakBefore/underline
# formatting.
$target.InsertXML($bulletsXml)
